$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "syntok" row (row 34) entirely - shifts rows 35-38 up to 34-37
$ws.Rows.Item(34).Delete()
